# Update the "Förändrad" (Changed) date column (C) for every data row
# (rows 2-116) on the active sheet from 2023-09-03 (serial 45172) to
# 2023-09-06 (serial 45175).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C116").Value = 45175
